# fix responce on select tab
$d = $word.ActiveDocument
$t = $d.Tables(1)

function Set-CellText($row, $text) {
    $cell = $t.Cell($row, 1)
    $r = $cell.Range
    $r.End = $r.End - 1
    $r.Text = $text
}

Set-CellText 1  "rfrfr frfrfr"
Set-CellText 2  "frfrf rfrfr"
Set-CellText 3  "3232 rfrfr"
Set-CellText 4  "cdcdcd"
Set-CellText 5  "cdcec@frfrf"
Set-CellText 6  "4343434"
Set-CellText 7  "rfrfrf"
Set-CellText 9  "___________________18-3-2020"
Set-CellText 11 "rfrfr frfrfr"
